# Remove column M from the alcohol measurement data: the whole column is
# deleted (not just cleared), so every column to its right (N) shifts one
# position to the left.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(13).Delete() | Out-Null

# Excel leaves the selection on the cell that now occupies the deleted
# column's position after a column delete.
$ws.Range("M1").Select() | Out-Null
